$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 90
$wsExhibit.Range("F3").Value = 376
$wsExhibit.Range("F4").Value = 4819
$wsExhibit.Range("F5").Value = 11
$wsExhibit.Range("F6").Value = 45

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 90
$wsAll.Range("F3").Value = 376
$wsAll.Range("F4").Value = 4819
$wsAll.Range("F6").Value = 11
$wsAll.Range("F8").Value = 45
